# Apply updated cryptocurrency price/volume figures scraped on
# Thu May 30 15:37:28 UTC 2024 (GitHub Actions data refresh).
#
# Column D ("Price") and column E ("Volume(1h)") are plain text cells
# (the sheet stores prices like "68.521.87" using dots as both group
# and decimal separators, so they can never be real Excel numbers).
# Where a new price happens to look like an ordinary decimal number
# (e.g. "1.00", "6.51") a leading apostrophe is used so Excel's
# COM layer stores it as literal text instead of silently coercing
# it to a number, matching the original text-only column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.521.87'
$ws.Range("E2").Value = '  +1.51%  '

$ws.Range("D3").Value = '3.779.76'
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = '''596.21'
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").Value = '''168.65'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").Value = '3.781.20'
$ws.Range("E7").Value = '  +0.70%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").Value = '''0.162'
$ws.Range("E10").Value = '  -1.96%  '

$ws.Range("D11").Value = '''6.51'
$ws.Range("E11").Value = '  +0.40%  '

$ws.Range("D12").Value = '''0.449'
$ws.Range("E12").Value = '  -1.44%  '

$ws.Range("E13").Value = '  -2.87%  '

$ws.Range("D14").Value = '''36.46'
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").Value = '4.408.59'
$ws.Range("E15").Value = '  -0.17%  '

$ws.Range("D16").Value = '3.768.78'
$ws.Range("E16").Value = '  -1.09%  '

$ws.Range("D17").Value = '68.444.50'
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").Value = '''18.24'
$ws.Range("E18").Value = '  -3.82%  '

$ws.Range("E19").Value = '  -2.71%  '

$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("D21").Value = '''10.94'
$ws.Range("E21").Value = '  +4.00%  '

$ws.Range("D22").Value = '''468.24'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("E23").Value = '  -2.78%  '

$ws.Range("D24").Value = '''84.93'
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25").Value = '''0.0000143'
$ws.Range("E25").Value = '  -4.13%  '

$ws.Range("D26").Value = '''2.23'
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("D27").Value = '''12.18'
$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("D28").Value = '''10.18'
$ws.Range("E28").Value = '  -1.17%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").Value = '3.922.32'
$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("E31").Value = '  -3.80%  '

$ws.Range("D32").Value = '''7.40'
$ws.Range("E32").Value = '  -2.98%  '

$ws.Range("E33").Value = '  -1.16%  '

$ws.Range("D34").Value = '''30.09'
$ws.Range("E34").Value = '  -0.92%  '

$ws.Range("D35").Value = '''9.26'
$ws.Range("E35").Value = '  +1.39%  '

$ws.Range("D36").Value = '''1.00'

$ws.Range("D37").Value = '3.727.40'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("E39").Value = '  -9.89%  '

$ws.Range("E40").Value = '  +1.26%  '

$ws.Range("D41").Value = '''1.01'
$ws.Range("E41").Value = '  +0.73%  '

$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("E43").Value = '  -0.34%  '

$ws.Range("E45").Value = '  -2.45%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '''1.97'
$ws.Range("E46").Value = '  +0.47%  '

$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").Value = '''43.69'
$ws.Range("E47").Value = '  +12.23%  '

$ws.Range("E48").Value = '  -1.49%  '

$ws.Range("D49").Value = '''407.27'
$ws.Range("E49").Value = '  +1.25%  '

$ws.Range("D50").Value = '''45.65'
$ws.Range("E50").Value = '  -1.33%  '

$ws.Range("D51").Value = '''145.82'
$ws.Range("E51").Value = '  +2.66%  '
